$d = $word.ActiveDocument

# 1. Publication year: 2009 -> 2008
$r1 = $d.Content.Find.Execute("2009", $true, $false, $false, $false, $false, $true, 1, $false, "2008", 2)
Write-Output "2009 -> 2008: $r1"

# 2. Paper title
$r2 = $d.Content.Find.Execute("Integrating Heterogeneous, Autonomous, Distributed Applications Using the DOM Prototype.", $true, $false, $false, $false, $false, $true, 1, $false, "A RISC Object Model for Object System Interoperation: Concepts and Applications.", 2)
Write-Output "title: $r2"

# 3. Author names
$r3 = $d.Content.Find.Execute("Mark F. Hornick, Joe D. Morrison, Farshad Nayeri", $true, $false, $false, $false, $false, $true, 1, $false, "Frank Manola, Sandra Heiler", 2)
Write-Output "authors: $r3"

# 4. Institution name
$r4 = $d.Content.Find.Execute("Alanus Hochschule für Kunst und Gesellschaft (Alfter)", $true, $false, $false, $false, $false, $true, 1, $false, "Technische Hochschule Aschaffenburg", 2)
Write-Output "institution: $r4"

# 5. Street address
$r5 = $d.Content.Find.Execute("Langenhorner Chaussee 86", $true, $false, $false, $false, $false, $true, 1, $false, "Gruenauer Strasse 59", 2)
Write-Output "address: $r5"
